# Apply the edits described in the commit: rename "unit_costs" sheet to
# "costs", drop the redundant "facilities" column header from both sheets,
# and apply a bold / centered / thin-bordered style to the header row and
# the facility-name column on both sheets.

$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("data")
$wsCosts = $wb.Worksheets.Item("unit_costs")

# 1. Rename the "unit_costs" sheet to "costs"
$wsCosts.Name = "costs"

# 2. Remove the "facilities" header label from A1 on both sheets - the
#    column no longer carries an explicit header.
$wsData.Cells.Item(1, 1).ClearContents()
$wsCosts.Cells.Item(1, 1).ClearContents()

# 3. Build the new header style once (bold, centered horizontally, top
#    vertically aligned, thin box border) on a single cell so the engine
#    folds all of the property changes into one cell style, then copy
#    that format onto the remaining header cells and the facility-name
#    column on each sheet.

# --- "data" sheet (columns B:K hold headers, A2:A9 hold facility names) ---
$styleSrc = $wsData.Range("B1")
$styleSrc.Font.Bold = $true
$styleSrc.Borders.LineStyle = 1
$styleSrc.HorizontalAlignment = -4108
$styleSrc.VerticalAlignment = -4160

$styleSrc.Copy()
$wsData.Range("C1:K1").PasteSpecial(-4122)
$styleSrc.Copy()
$wsData.Range("A2:A9").PasteSpecial(-4122)

# --- "costs" sheet (columns B:I hold headers, A2:A9 hold facility names) ---
# Reuse the style already built on $styleSrc (copy/paste its format) rather
# than re-deriving it through a fresh sequence of property assignments, so
# the workbook ends up with a single shared cell style instead of leftover
# intermediate ones.
$styleSrc.Copy()
$wsCosts.Range("B1:I1").PasteSpecial(-4122)
$styleSrc.Copy()
$wsCosts.Range("A2:A9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 4. Update the selected cell on each sheet.
$wsData.Activate()
$wsData.Range("E19").Select()

$wsCosts.Activate()
$wsCosts.Range("I21").Select()

# Leave "data" as the active (tab-selected) sheet, matching the original
# workbook where the first sheet is the one shown on open.
$wsData.Activate()
